# Second commit: populate Sheet1!D10 with text and move the active selection to F16.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D10").Value = "shreyansh preparation"
$ws.Range("F16").Select()
